$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SARIF2003 (ProvideVersionControlProvenance, currently row 44) now has its
# rule messages written, so "Message code status" (column H) moves from
# TODO to DONE. Match the existing DONE formatting (green fill) used
# elsewhere in the column by pasting the format from a DONE cell in the
# same row (column F) instead of hand-rolling font/fill properties.
$ws.Range("H44").Value = "DONE"
$ws.Range("F44").Copy() | Out-Null
$ws.Range("H44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The SARIF2xxx rules block used to start with three blank spacer rows
# (37:39) above it; trim it down to the single blank spacer row used
# everywhere else on the sheet by removing two of them. This shifts every
# row from 40:73 up to 38:71.
$ws.Rows("37:38").Delete() | Out-Null

# Reflect where the edit was made.
$ws.Range("H42").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
